# The workbook gains one new data row. A new row is inserted at row 594
# (pushing the existing rows 594-630 down to 595-631) and is populated with
# a new "Coliflor" price observation for "Provincia de Melipilla".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 594, shifting everything
# from 594 downward by one row (old row 594 becomes row 595, ..., old row
# 630 becomes row 631). This also grows the sheet dimension to A1:R631.
$ws.Rows.Item(594).Insert()

# Populate the newly inserted row 594 with the new record's values.
$ws.Cells.Item(594, 1).Value  = 9
$ws.Cells.Item(594, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(594, 3).Value  = "Metropolitana"
$ws.Cells.Item(594, 4).Value  = 44706
$ws.Cells.Item(594, 5).Value  = 13
$ws.Cells.Item(594, 6).Value  = 100112008
$ws.Cells.Item(594, 7).Value  = "Coliflor"
$ws.Cells.Item(594, 8).Value  = "Sin especificar"
$ws.Cells.Item(594, 9).Value  = "Primera"
$ws.Cells.Item(594, 10).Value = 2600
$ws.Cells.Item(594, 11).Value = 1000
$ws.Cells.Item(594, 12).Value = 1100
$ws.Cells.Item(594, 13).Value = 1054
$ws.Cells.Item(594, 14).Value = "$/unidad"
$ws.Cells.Item(594, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(594, 16).Value = 1054
$ws.Cells.Item(594, 17).Value = 1
$ws.Cells.Item(594, 18).Value = "Hortaliza"
